# Patch release 3.0.1 edits for the PACT Simplified Data Model workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bump version number in A2: 3.0.0 -> 3.0.1
$ws.Range("A2").Value = "3.0.1"

# 2. Clarify the `declaredUnitOfMeasurement` description in E21.
$e21 = 'The unit of measurement of the product. Together with `declaredUnitAmount` this defines the ''declared unit'' of the product. Emissions in this carbon footprint are expressed in kgCO2e per ''declared unit''. ' + "`n" + 'For example: a PCF for a 12.5 liter bottle of Ethanol states 2 kg of CO2e in emissions. In this case the declared unit is 12.5 liter Ethanol, thus the `declaredUnitOfMeasurement` is "liter", and the `declaredUnitAmount` is "12.5". The `pcfIncludingBiogenicUptake` is "2.0" kgCO2e per "12.5 liter" of Ethanol.'
$ws.Range("E21").Value = $e21

# 3. Remove invalid geography example values (region/sub-region, country, country subdivision).
$ws.Range("J26:L26").Value = ""
$ws.Range("J27").Value = ""
$ws.Range("L27:M27").Value = ""
$ws.Range("K28").Value = ""
$ws.Range("M28").Value = ""

# 4. Extend the accepted cross-sectoral standards list (I57) and example (M57).
$ws.Range("I57").Value = "array: ISO14067|ISO14083|ISO14040-44|GHGP-Product|PEF|PACT-1.0|PACT-2.0|PACT-3.0|PAS2050|..." + "`n" + "(string)"
$ws.Range("M57").Value = "['PEF', 'FUTURE-STANDARD']"
